# Publish the testing slide (slide 9, "Testing") with the body copy that
# explains the team's testing philosophy / approach.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$shp = $s.Shapes.Item(2)

# The placeholder currently holds a single empty paragraph
# (<a:p><a:endParaRPr lang="en-US" dirty="0"/></a:p>). Using InsertAfter
# (instead of a plain TextRange.Text assignment) keeps that empty
# paragraph's trailing endParaRPr intact after the new content instead of
# it being clobbered.
$tr = $shp.TextFrame.TextRange
[void]$tr.InsertAfter("Our philosophy to testing has been to focus on getting working code first and test as much as possible during the coding process.")

$tr = $shp.TextFrame.TextRange
[void]$tr.InsertAfter("`rBecause this is our first time using Junit we elected to test manually up until final debugging. ")

$tr = $shp.TextFrame.TextRange
[void]$tr.InsertAfter("`rWe have regularly used manual testing to fix collisions and other errors that have come up during the code writing phase of our project. ")

$tr = $shp.TextFrame.TextRange
[void]$tr.InsertAfter("`rWe are working diligently to have Junit testing completed ")

# Appended as its own trailing run (the source has two <a:r> runs in this
# final paragraph) rather than being merged into the previous insert.
$tr = $shp.TextFrame.TextRange
[void]$tr.InsertAfter("before the final due date.")
